# Updates the 'Pais' COVID-19 country/provincia dataset:
#  - Refreshes the 'Datos actualizados' timestamp in A1
#  - Re-sorts several countries (rows 81-180) into their correct
#    ranking position, swapping each row's name + stats so the
#    table stays ordered by case counts (per the source diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh 'last updated' timestamp
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 20 de Marzo de 2020 a las 02:16'

# Swap Albania / Bosnia y Herzegovina (rows 81-82)
$ws.Cells.Item(81, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(81, 3).Value = 25
$ws.Cells.Item(81, 4).Value = 2
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(82, 1).Value = 'Albania'
$ws.Cells.Item(82, 3).Value = 5
$ws.Cells.Item(82, 4).Value = 0
$ws.Cells.Item(82, 6).Value = 2
$ws.Cells.Item(82, 8).Value = 2

# Insert Paraguay before Montenegro/Maldivas/Camerun, shifting them down (rows 115-118)
$ws.Cells.Item(115, 1).Value = 'Paraguay'
$ws.Cells.Item(115, 3).Value = 2
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(116, 1).Value = 'Montenegro'
$ws.Cells.Item(116, 3).Value = 5
$ws.Cells.Item(117, 1).Value = 'Maldivas'
$ws.Cells.Item(118, 1).Value = 'Camerun'
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 6).Value = 0

# Swap Guam / Honduras (rows 119-120)
$ws.Cells.Item(119, 1).Value = 'Guam'
$ws.Cells.Item(119, 3).Value = 4
$ws.Cells.Item(120, 1).Value = 'Honduras'
$ws.Cells.Item(120, 3).Value = 3

# Move Kenia up before Mauricio/Etiopia, then re-rank the
# Polinesia Francesa/Tanzania/Guinea Ecuatorial/Seychelles/Puerto Rico/Mongolia block (rows 130-138)
$ws.Cells.Item(130, 1).Value = 'Kenia'
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(131, 1).Value = 'Mauricio'
$ws.Cells.Item(131, 3).Value = 4
$ws.Cells.Item(132, 1).Value = 'Etiopia'
$ws.Cells.Item(132, 3).Value = 1
$ws.Cells.Item(133, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(134, 1).Value = 'Tanzania'
$ws.Cells.Item(134, 3).Value = 3
$ws.Cells.Item(135, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(135, 3).Value = 2
$ws.Cells.Item(136, 1).Value = 'Seychelles'
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(137, 1).Value = 'Puerto Rico'
$ws.Cells.Item(138, 1).Value = 'Mongolia'
$ws.Cells.Item(138, 3).Value = 0

# Swap Aruba / Guyana (rows 140-141)
$ws.Cells.Item(140, 1).Value = 'Aruba'
$ws.Cells.Item(140, 4).Value = 1
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(141, 1).Value = 'Guyana'
$ws.Cells.Item(141, 4).Value = 0
$ws.Cells.Item(141, 8).Value = 1

# Re-rank the Bahamas/San Martin (Francesa)/San Bartolome/Kirguistan/
# Islas Virgenes EEUU/Namibia/Congo/Gabon block (rows 143-150)
$ws.Cells.Item(143, 1).Value = 'Bahamas'
$ws.Cells.Item(143, 3).Value = 2
$ws.Cells.Item(144, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(145, 1).Value = 'San Bartolome'
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(146, 1).Value = 'Kirguistan'
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(149, 1).Value = 'Congo'
$ws.Cells.Item(149, 3).Value = 2
$ws.Cells.Item(150, 1).Value = 'Gabon'

# Swap Islas Caimanes / Curazao (rows 151-152)
$ws.Cells.Item(151, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(152, 1).Value = 'Curazao'
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 7).Value = 1

# Re-rank the Benin/Mauritania/Liberia/Zambia/Groenlandia/Santa Lucia/
# Nueva Caledonia/Bermudas block (rows 153-160)
$ws.Cells.Item(153, 1).Value = 'Benin'
$ws.Cells.Item(154, 1).Value = 'Mauritania'
$ws.Cells.Item(155, 1).Value = 'Liberia'
$ws.Cells.Item(156, 1).Value = 'Zambia'
$ws.Cells.Item(157, 1).Value = 'Groenlandia'
$ws.Cells.Item(158, 1).Value = 'Santa Lucia'
$ws.Cells.Item(159, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(160, 1).Value = 'Bermudas'

# Re-rank the Nicaragua/Fiyi/Republica del Chad/San Vicente y las Granadinas/
# Republica de Yibuti/Santa Sede/Guinea/Niger/Surinam block (rows 162-170)
$ws.Cells.Item(162, 1).Value = 'Nicaragua'
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(163, 1).Value = 'Fiyi'
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(164, 1).Value = 'Republica del Chad'
$ws.Cells.Item(164, 3).Value = 1
$ws.Cells.Item(165, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(166, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(168, 1).Value = 'Guinea'
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(169, 1).Value = 'Niger'
$ws.Cells.Item(170, 1).Value = 'Surinam'

# Re-rank the Montserrat/Gambia/San Martin (Holandesa)/Togo/El Salvador/
# Antigua y Barbuda/Somalia/Isla de Man/Butan block (rows 172-180)
$ws.Cells.Item(172, 1).Value = 'Montserrat'
$ws.Cells.Item(173, 1).Value = 'Gambia'
$ws.Cells.Item(174, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(175, 1).Value = 'Togo'
$ws.Cells.Item(176, 1).Value = 'El Salvador'
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(177, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(178, 1).Value = 'Somalia'
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(179, 1).Value = 'Isla de Man'
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(180, 1).Value = 'Butan'
$ws.Cells.Item(180, 3).Value = 0
